# PV-261: Correct PV-Test-03 so it has a valid sheet name, and make sure
# the worksheet containing the plan data (not the "Dummy" sheet) is the
# active/selected sheet when the workbook is opened.

$wb = $excel.ActiveWorkbook

# The data sheet was incorrectly named "PV-Test-01" - rename it to match
# the file name "PV-Test-03".
$dataSheet = $wb.Worksheets.Item("PV-Test-01")
$dataSheet.Name = "PV-Test-03"

# Make the data sheet the active sheet (previously "Dummy" was active),
# so a naive "read the active sheet" import picks up the right data.
$dataSheet.Activate()
